$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.173.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'1.834.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D5").Value = "'241.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'0.6606"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.07410"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "'0.2936"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").Value = "'1.824.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "'4.981"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'0.6659"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "'82.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "'6.098"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "'0.000008511"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").Value = "'29.168.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'2.092.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'227.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'12.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'7.081"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'159.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'8.602"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.1399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'1.514"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'4.111"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D32").Value = "'1.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "'0.05259"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "'1.863"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").Value = "'0.7353"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "'1.144"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").Value = "'2.660"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").Value = "'1.299.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "'0.01791"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").Value = "'0.9202"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'6.060"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").Value = "'0.08434"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.02%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'102.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "'1.992.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'0.5139"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "'0.00000000121"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.753"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'63.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
